$wb = $excel.ActiveWorkbook

# --- Reference to the existing "总计" summary sheet ---
$total = $wb.Worksheets.Item(1)   # "总计"

# --- 1. Update the "总计" summary sheet ---
# Push the existing 2022-Q2 row (row 2) down to row 3, keeping its formatting,
# then write the new 2022-Q3 figures into row 2.
$total.Range("A2:D2").Copy($total.Range("A3:D3"))
$total.Range("A3").Value = 1

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 6
$total.Range("D2").Value = 2.09

# --- 2. Insert the new "2022-Q3" worksheet between "总计" and "2022-Q2" ---
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# Copy header formatting (bold + border style used on the "总计" sheet) onto
# the new sheet's header row and first data column.
$total.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$q3.Range("A2:A7").PasteSpecial(-4122)

# Header row
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Columns B-G hold text (fund codes with leading zeros, numbers formatted as
# text) in the source data, so force text formatting before assigning values
# to avoid Excel auto-converting them to numbers.
$q3.Range("B2:G7").NumberFormat = "@"

# Data rows
$q3Data = @(
    @(0, "006624", "中泰玉衡价值优选混合A", "23.07", "88.61", "3.41", "0.7867", 10),
    @(1, "013776", "中泰兴为价值精选混合A", "12.11", "90.90", "4.90", "0.5934", 4),
    @(2, "010728", "中泰兴诚价值一年持有期混合A", "6.63", "89.26", "4.99", "0.3308", 10),
    @(3, "013777", "中泰兴为价值精选混合C", "5.78", "90.90", "4.90", "0.2832", 4),
    @(4, "010729", "中泰兴诚价值一年持有期混合C", "1.27", "89.26", "4.99", "0.0634", 10),
    @(5, "016090", "中泰玉衡价值优选混合C", "1.01", "88.61", "3.41", "0.0344", 10)
)

$r = 2
foreach ($row in $q3Data) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
